$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order/content for rows 2-9 (B = id, C = speaker_variant)
$data = @(
    @{B="#koor"; C="Koor"},
    @{B="#eerste-draager"; C="Eerste Draager"},
    @{B="#de-heerschzugtige"; C="De Heerschzugtige"},
    @{B="#de-geldzugtige"; C="De Geldzugtige"},
    @{B="#de-minzieke-juffer"; C="De Minzieke Juffer"},
    @{B="#twede-draager"; C="Twede Draager"},
    @{B="#de-minzieke-juffer,-de-heerschzugtige"; C="De Minzieke Juffer, De Heerschzugtige"},
    @{B="#de-kwakzalver"; C="De Kwakzalver"}
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i].B
    $ws.Cells.Item($row, 3).Value = $data[$i].C
    $ws.Cells.Item($row, 4).Value = ""
}
